{"js": "// Replace the date and each three-digit-by-one-digit multiplication\n// problem in the worksheet with the updated values from the commit.\nconst replacements = [\n  [\"2025-08-12 Tuesday\", \"2025-08-13 Wednesday\"],\n  [\"357\\u00d74=\", \"101\\u00d73=\"],\n  [\"562\\u00d77=\", \"589\\u00d72=\"],\n  [\"961\\u00d74=\", \"172\\u00d73=\"],\n  [\"981\\u00d72=\", \"119\\u00d79=\"],\n  [\"811\\u00d74=\", \"108\\u00d77=\"],\n  [\"983\\u00d77=\", \"782\\u00d75=\"],\n  [\"486\\u00d77=\", \"220\\u00d75=\"],\n  [\"658\\u00d76=\", \"233\\u00d79=\"],\n  [\"846\\u00d73=\", \"568\\u00d73=\"],\n  [\"949\\u00d77=\", \"309\\u00d74=\"],\n  [\"610\\u00d77=\", \"145\\u00d76=\"],\n  [\"376\\u00d72=\", \"851\\u00d77=\"],\n  [\"763\\u00d78=\", \"950\\u00d72=\"],\n  [\"757\\u00d78=\", \"718\\u00d76=\"],\n  [\"734\\u00d76=\", \"737\\u00d75=\"],\n  [\"105\\u00d75=\", \"597\\u00d73=\"],\n  [\"253\\u00d72=\", \"495\\u00d74=\"],\n  [\"231\\u00d76=\", \"741\\u00d72=\"],\n  [\"739\\u00d72=\", \"852\\u00d78=\"],\n  [\"124\\u00d75=\", \"440\\u00d79=\"],\n  [\"462\\u00d78=\", \"372\\u00d74=\"],\n  [\"538\\u00d75=\", \"177\\u00d72=\"],\n  [\"412\\u00d76=\", \"289\\u00d75=\"],\n  [\"737\\u00d73=\", \"812\\u00d77=\"],\n  [\"134\\u00d75=\", \"431\\u00d76=\"],\n];\n\nfor (const [from, to] of replacements) {\n  const results = context.document.body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and each three-digit-by-one-digit\n# multiplication problem to the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-12 Tuesday\", \"2025-08-13 Wednesday\"),\n    @(\"357\u00d74=\", \"101\u00d73=\"),\n    @(\"562\u00d77=\", \"589\u00d72=\"),\n    @(\"961\u00d74=\", \"172\u00d73=\"),\n    @(\"981\u00d72=\", \"119\u00d79=\"),\n    @(\"811\u00d74=\", \"108\u00d77=\"),\n    @(\"983\u00d77=\", \"782\u00d75=\"),\n    @(\"486\u00d77=\", \"220\u00d75=\"),\n    @(\"658\u00d76=\", \"233\u00d79=\"),\n    @(\"846\u00d73=\", \"568\u00d73=\"),\n    @(\"949\u00d77=\", \"309\u00d74=\"),\n    @(\"610\u00d77=\", \"145\u00d76=\"),\n    @(\"376\u00d72=\", \"851\u00d77=\"),\n    @(\"763\u00d78=\", \"950\u00d72=\"),\n    @(\"757\u00d78=\", \"718\u00d76=\"),\n    @(\"734\u00d76=\", \"737\u00d75=\"),\n    @(\"105\u00d75=\", \"597\u00d73=\"),\n    @(\"253\u00d72=\", \"495\u00d74=\"),\n    @(\"231\u00d76=\", \"741\u00d72=\"),\n    @(\"739\u00d72=\", \"852\u00d78=\"),\n    @(\"124\u00d75=\", \"440\u00d79=\"),\n    @(\"462\u00d78=\", \"372\u00d74=\"),\n    @(\"538\u00d75=\", \"177\u00d72=\"),\n    @(\"412\u00d76=\", \"289\u00d75=\"),\n    @(\"737\u00d73=\", \"812\u00d77=\"),\n    @(\"134\u00d75=\", \"431\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
